$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4 (ALC)
$ws.Range("H4").Value = 968.65216
$ws.Range("I4").Value = 968.65216
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 968.65216
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -854.65216
$ws.Range("N4").ClearContents()

# Row 62 (ALC)
$ws.Range("H62").Value = 1592.5
$ws.Range("I62").Value = 1238.75
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 1238.75
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -614.75
$ws.Range("N62").Value = -3548

# Row 65 (ALC)
$ws.Range("H65").Value = 1592.5
$ws.Range("I65").Value = 1238.75
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 6193.75
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -3073.75
$ws.Range("N65").Value = -17740

# Row 129 (ALC)
$ws.Range("H129").Value = 3084.6
$ws.Range("I129").Value = 10535
$ws.Range("J129").Value = 955.9143
$ws.Range("K129").Value = 31605
$ws.Range("L129").Value = 2867.7429
$ws.Range("M129").Value = -26605
$ws.Range("N129").Value = -12867.7429

# Row 138 (ALC)
$ws.Range("H138").Value = 2127.6064
$ws.Range("I138").Value = 1055.0731
$ws.Range("J138").Value = 2957.302
$ws.Range("K138").Value = 3165.2193
$ws.Range("L138").Value = 8871.906000000001
$ws.Range("M138").Value = 1974.7807
$ws.Range("N138").Value = -19151.906

# Row 140 (ALC)
$ws.Range("H140").Value = 47637.5
$ws.Range("J140").Value = 47637.5
$ws.Range("L140").Value = 47637.5
$ws.Range("N140").Value = -57997.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 25881.242
$ws.Range("I32").Value = 4857.7666
$ws.Range("J32").Value = 152022.1
$ws.Range("K32").Value = 4857.7666
$ws.Range("L32").Value = 152022.1
$ws.Range("M32").Value = -4570.7666
$ws.Range("N32").Value = -152596.1

# Row 61 (ARM)
$ws.Range("H61").Value = 1322.0834
$ws.Range("I61").Value = 950.1111
$ws.Range("J61").Value = 2438
$ws.Range("K61").Value = 950.1111
$ws.Range("L61").Value = 2438
$ws.Range("M61").Value = -738.1111
$ws.Range("N61").Value = -2862

# Row 136 (ARM)
$ws.Range("H136").Value = 1322.0834
$ws.Range("I136").Value = 950.1111
$ws.Range("J136").Value = 2438
$ws.Range("K136").Value = 2850.3333
$ws.Range("L136").Value = 7314
$ws.Range("M136").Value = -300.3332999999998
$ws.Range("N136").Value = -12414

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (BSM)
$ws.Range("H22").Value = 550
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -327
$ws.Range("N22").Value = -946

# Row 99 (BSM)
$ws.Range("H99").Value = 1642.75
$ws.Range("I99").Value = 1440.3334
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 1440.3334
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = 57.66660000000002
$ws.Range("N99").Value = -5246

# Row 107 (BSM)
$ws.Range("H107").Value = 90910250
$ws.Range("I107").Value = 200001180
$ws.Range("J107").Value = 1138.6666
$ws.Range("K107").Value = 200001180
$ws.Range("L107").Value = 1138.6666
$ws.Range("M107").Value = -199999260
$ws.Range("N107").Value = -4978.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 28216.932
$ws.Range("I31").Value = 1090.9678
$ws.Range("J31").Value = 58249.25
$ws.Range("K31").Value = 1090.9678
$ws.Range("L31").Value = 58249.25
$ws.Range("M31").Value = -795.9677999999999
$ws.Range("N31").Value = -58839.25

# Row 34 (CRP)
$ws.Range("H34").Value = 28216.932
$ws.Range("I34").Value = 1090.9678
$ws.Range("J34").Value = 58249.25
$ws.Range("K34").Value = 1090.9678
$ws.Range("L34").Value = 58249.25
$ws.Range("M34").Value = -888.9677999999999
$ws.Range("N34").Value = -58653.25

# Row 58 (CRP)
$ws.Range("H58").Value = 4563.2607
$ws.Range("I58").Value = 1223.2858
$ws.Range("K58").Value = 1223.2858
$ws.Range("M58").Value = -1020.2858

# Row 62 (CRP)
$ws.Range("H62").Value = 2454.7778
$ws.Range("I62").Value = 2098.3333
$ws.Range("J62").Value = 2633
$ws.Range("K62").Value = 2098.3333
$ws.Range("L62").Value = 2633
$ws.Range("M62").Value = -1474.3333
$ws.Range("N62").Value = -3881

# Row 65 (CRP)
$ws.Range("H65").Value = 2454.7778
$ws.Range("I65").Value = 2098.3333
$ws.Range("J65").Value = 2633
$ws.Range("K65").Value = 10491.6665
$ws.Range("L65").Value = 13165
$ws.Range("M65").Value = -7371.666499999999
$ws.Range("N65").Value = -19405

# Row 99 (CRP)
$ws.Range("H99").Value = 11964.182
$ws.Range("I99").Value = 3230.6667
$ws.Range("J99").Value = 15239.25
$ws.Range("K99").Value = 3230.6667
$ws.Range("L99").Value = 15239.25
$ws.Range("M99").Value = -1732.6667
$ws.Range("N99").Value = -18235.25

# Row 122 (CRP)
$ws.Range("H122").Value = 944.4
$ws.Range("I122").Value = 688.8
$ws.Range("K122").Value = 2066.4
$ws.Range("M122").Value = 383.6000000000004

# Row 126 (CRP)
$ws.Range("H126").Value = 11964.182
$ws.Range("I126").Value = 3230.6667
$ws.Range("J126").Value = 15239.25
$ws.Range("K126").Value = 9692.000100000001
$ws.Range("L126").Value = 45717.75
$ws.Range("M126").Value = -7222.000100000001
$ws.Range("N126").Value = -50657.75

# Row 132 (CRP)
$ws.Range("H132").Value = 2423.1052
$ws.Range("I132").Value = 2018
$ws.Range("J132").Value = 3417.4546
$ws.Range("K132").Value = 6054
$ws.Range("L132").Value = 10252.3638
$ws.Range("M132").Value = -3524
$ws.Range("N132").Value = -15312.3638

# Row 136 (CRP)
$ws.Range("H136").Value = 4563.2607
$ws.Range("I136").Value = 1223.2858
$ws.Range("K136").Value = 3669.8574
$ws.Range("M136").Value = -1119.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 126 (GSM)
$ws.Range("H126").Value = 1636491.6
$ws.Range("I126").Value = 2057
$ws.Range("J126").Value = 4204889
$ws.Range("K126").Value = 6171
$ws.Range("L126").Value = 12614667
$ws.Range("M126").Value = -3701
$ws.Range("N126").Value = -12619607

# Row 132 (GSM)
$ws.Range("H132").Value = 1193.0741
$ws.Range("I132").Value = 1019.2273
$ws.Range("J132").Value = 1958
$ws.Range("K132").Value = 3057.6819
$ws.Range("L132").Value = 5874
$ws.Range("M132").Value = -527.6819
$ws.Range("N132").Value = -10934

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (LTW)
$ws.Range("H40").Value = 73450
$ws.Range("I40").Value = 251125
$ws.Range("J40").Value = 2380
$ws.Range("K40").Value = 251125
$ws.Range("L40").Value = 2380
$ws.Range("M40").Value = -250989
$ws.Range("N40").Value = -2652

# Row 132 (LTW)
$ws.Range("H132").Value = 1675.8914
$ws.Range("I132").Value = 1519.9722
$ws.Range("J132").Value = 2237.2
$ws.Range("K132").Value = 4559.9166
$ws.Range("L132").Value = 6711.599999999999
$ws.Range("M132").Value = -2029.9166
$ws.Range("N132").Value = -11771.6

# Row 136 (LTW)
$ws.Range("H136").Value = 1560.6
$ws.Range("I136").Value = 1457.6364
$ws.Range("J136").Value = 2046
$ws.Range("K136").Value = 4372.9092
$ws.Range("L136").Value = 6138
$ws.Range("M136").Value = -1822.9092
$ws.Range("N136").Value = -11238

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 2400.6924
$ws.Range("I122").Value = 1400.2858
$ws.Range("J122").Value = 3567.8333
$ws.Range("K122").Value = 4200.857400000001
$ws.Range("L122").Value = 10703.4999
$ws.Range("M122").Value = -1750.857400000001
$ws.Range("N122").Value = -15603.4999

# Row 126 (WVR)
$ws.Range("H126").Value = 2372.375
$ws.Range("I126").Value = 2372.375
$ws.Range("K126").Value = 7117.125
$ws.Range("M126").Value = -4647.125

Write-Host "Applied scheduled market-price refresh to 29 rows across 7 sheets"
